# Weekly price-report update: insert 4 new records (week of D=45013, i.e.
# 2023-03-28) for "Terminal Hortofrutícola Agro Chillán" / Manzana at the
# top of the existing data block (row 1088), pushing the historical rows
# down by 4 (old 1088..1174 -> new 1092..1178).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 1088; this shifts the existing
# rows 1088-1174 down to 1092-1178 and grows the used range accordingly.
$ws.Rows.Item(1088).Resize(4).Insert()

# Column layout (fixed across this whole sheet/subset):
# A Mercado ID | B Mercado | C Región | D Fecha | E Codreg | F Tipo
# G Producto ID | H Producto | I Categoría ID | J Categoría | K Variedad
# L Calidad | M Volumen | N Precio mínimo | O Precio máximo
# P Precio promedio ponderado | Q Unidad de comercialización | R Origen
# S Precio $/Kg | T Kg / unidad

$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$tipo        = "Fruta"
$productoId  = 100104
$producto    = "Frutos de pepita"
$categoriaId = 100104002
$categoria   = "Manzana"
$unidad      = "$/caja 16 kilos empedrada"
$kgUnidad    = 16

$newRows = @(
    @{ Row=1088; Fecha=45013; Variedad="Granny Smith"; Calidad="Especial"; Volumen=50;  PMin=12000; PMax=12000; PProm=12000; Origen="Región de O'Higgins"; PrecioKg=750 },
    @{ Row=1089; Fecha=45013; Variedad="Granny Smith"; Calidad="Primera";  Volumen=100; PMin=10000; PMax=10000; PProm=10000; Origen="Región de O'Higgins"; PrecioKg=625 },
    @{ Row=1090; Fecha=45013; Variedad="Royal Gala";   Calidad="Especial"; Volumen=50;  PMin=12000; PMax=12000; PProm=12000; Origen="Región de O'Higgins"; PrecioKg=750 },
    @{ Row=1091; Fecha=45013; Variedad="Royal Gala";   Calidad="Primera";  Volumen=100; PMin=10000; PMax=10000; PProm=10000; Origen="Región de O'Higgins"; PrecioKg=625 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
